# Refresh captured Universalis market-price snapshots (and the leve-profit
# calculations derived from them) for each job sheet's Leve table.
#
# Columns (per Table_<JOB>, A1:N141):
#   H currentAveragePrice    I currentAveragePriceNQ  J currentAveragePriceHQ
#   K LevePriceNQ            L LevePriceHQ            M LeveProfitNQ
#   N LeveProfitHQ
# Cells that have no value for a given leve are left/made blank (ClearContents),
# matching how the upstream scraper omits zero/not-applicable data points.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 566
$ws.Range("I5").Value = 130
$ws.Range("J5").Value = 1002
$ws.Range("K5").Value = 130
$ws.Range("L5").Value = 1002
$ws.Range("M5").Value = -15
$ws.Range("N5").Value = -1232
# Row 32
$ws.Range("H32").Value = 3967.4
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3967.4
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3967.4
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4619.4
# Row 38
$ws.Range("H38").Value = 2915
$ws.Range("I38").Value = 2377.6924
$ws.Range("K38").Value = 7133.0772
$ws.Range("M38").Value = -6761.0772
# Row 39
$ws.Range("H39").Value = 1159.9
$ws.Range("I39").Value = 496.6
$ws.Range("J39").Value = 1823.2
$ws.Range("K39").Value = 1489.8
$ws.Range("L39").Value = 5469.6
$ws.Range("M39").Value = -1193.8
$ws.Range("N39").Value = -6061.6
# Row 80
$ws.Range("H80").Value = 1825.1
$ws.Range("I80").Value = 858.6667
$ws.Range("J80").Value = 2239.2856
$ws.Range("K80").Value = 2576.0001
$ws.Range("L80").Value = 6717.8568
$ws.Range("M80").Value = -1578.0001
$ws.Range("N80").Value = -8713.856800000001
# Row 83
$ws.Range("H83").Value = 1825.1
$ws.Range("I83").Value = 858.6667
$ws.Range("J83").Value = 2239.2856
$ws.Range("K83").Value = 7728.0003
$ws.Range("L83").Value = 20153.5704
$ws.Range("M83").Value = -2736.0003
$ws.Range("N83").Value = -30137.5704
# Row 92
$ws.Range("H92").Value = 707.1667
$ws.Range("I92").Value = 534.9286
$ws.Range("K92").Value = 534.9286
$ws.Range("M92").Value = 713.0714
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
# Row 135
$ws.Range("H135").Value = 1949.7646
$ws.Range("I135").Value = 2083.2856
$ws.Range("J135").Value = 1326.6666
$ws.Range("K135").Value = 18749.5704
$ws.Range("L135").Value = 11939.9994
$ws.Range("M135").Value = -16214.5704
$ws.Range("N135").Value = -17009.9994
# Row 137
$ws.Range("H137").Value = 1449.1666
$ws.Range("I137").Value = 1314.5714
$ws.Range("J137").Value = 1637.6
$ws.Range("K137").Value = 3943.7142
$ws.Range("L137").Value = 4912.799999999999
$ws.Range("M137").Value = -1393.7142
$ws.Range("N137").Value = -10012.8

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3921.9167
$ws.Range("I45").Value = 3921.9167
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3921.9167
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3544.9167
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
# Row 94
$ws.Range("H94").Value = 1395.0714
$ws.Range("I94").Value = 1228.3334
$ws.Range("J94").Value = 1695.2
$ws.Range("K94").Value = 1228.3334
$ws.Range("L94").Value = 1695.2
$ws.Range("M94").Value = -777.3334
$ws.Range("N94").Value = -2597.2
# Row 95
$ws.Range("H95").Value = 16257.5
$ws.Range("J95").Value = 16257.5
$ws.Range("L95").Value = 16257.5
$ws.Range("N95").Value = -21749.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 289.45456
$ws.Range("I7").Value = 128.2
$ws.Range("K7").Value = 128.2
$ws.Range("M7").Value = -15.19999999999999
# Row 81
$ws.Range("H81").Value = 78475
$ws.Range("J81").Value = 78475
$ws.Range("L81").Value = 78475
$ws.Range("N81").Value = -80471
# Row 84
$ws.Range("H84").Value = 78475
$ws.Range("J84").Value = 78475
$ws.Range("L84").Value = 235425
$ws.Range("N84").Value = -245409
# Row 88
$ws.Range("H88").Value = 18146
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 18146
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 18146
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -18958
# Row 91
$ws.Range("H91").Value = 18146
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 18146
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 18146
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -20954
# Row 92
$ws.Range("H92").Value = 29996.5
$ws.Range("J92").Value = 29996.5
$ws.Range("L92").Value = 29996.5
$ws.Range("N92").Value = -34988.5
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 99
$ws.Range("H99").Value = 734.6667
$ws.Range("I99").Value = 695.53845
$ws.Range("J99").Value = 836.4
$ws.Range("K99").Value = 695.53845
$ws.Range("L99").Value = 836.4
$ws.Range("M99").Value = 802.46155
$ws.Range("N99").Value = -3832.4
# Row 126
$ws.Range("H126").Value = 734.6667
$ws.Range("I126").Value = 695.53845
$ws.Range("J126").Value = 836.4
$ws.Range("K126").Value = 2086.61535
$ws.Range("L126").Value = 2509.2
$ws.Range("M126").Value = 383.38465
$ws.Range("N126").Value = -7449.2

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 192
$ws.Range("I12").Value = 250
$ws.Range("K12").Value = 750
$ws.Range("M12").Value = -577
# Row 57
$ws.Range("H57").Value = 17857.428
$ws.Range("I57").Value = 15000.667
$ws.Range("K57").Value = 45002.001
$ws.Range("M57").Value = -44443.001
# Row 68
$ws.Range("H68").Value = 3742.9333
$ws.Range("J68").Value = 3742.9333
$ws.Range("L68").Value = 11228.7999
$ws.Range("N68").Value = -12850.7999
# Row 71
$ws.Range("H71").Value = 3742.9333
$ws.Range("J71").Value = 3742.9333
$ws.Range("L71").Value = 33686.3997
$ws.Range("N71").Value = -41798.3997
# Row 87
$ws.Range("H87").Value = 13995
$ws.Range("I87").Value = 13995
$ws.Range("K87").Value = 41985
$ws.Range("M87").Value = -40737
# Row 90
$ws.Range("H90").Value = 13995
$ws.Range("I90").Value = 13995
$ws.Range("K90").Value = 125955
$ws.Range("M90").Value = -119715
# Row 113
$ws.Range("H113").Value = 1156.3077
$ws.Range("J113").Value = 1169.3334
$ws.Range("L113").Value = 3508.0002
$ws.Range("N113").Value = -7848.0002
# Row 137
$ws.Range("H137").Value = 2810.2856
$ws.Range("I137").Value = 2392.5
$ws.Range("K137").Value = 7177.5
$ws.Range("M137").Value = -2077.5

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 20200.8
$ws.Range("I93").Value = 20000
$ws.Range("K93").Value = 20000
$ws.Range("M93").Value = -18128
# Row 132
$ws.Range("H132").Value = 4273.316
$ws.Range("I132").Value = 4540.7646
$ws.Range("K132").Value = 13622.2938
$ws.Range("M132").Value = -11092.2938
# Row 141
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 1829.2858
$ws.Range("I9").Value = 1313.75
$ws.Range("J9").Value = 2516.6667
$ws.Range("K9").Value = 1313.75
$ws.Range("L9").Value = 2516.6667
$ws.Range("M9").Value = -1089.75
$ws.Range("N9").Value = -2964.6667
# Row 12
$ws.Range("H12").Value = 1200
$ws.Range("J12").Value = 1500
$ws.Range("L12").Value = 1500
$ws.Range("N12").Value = -1840
# Row 19
$ws.Range("H19").Value = 1745
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1745
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1745
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -2085
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 31
$ws.Range("H31").Value = 3999
$ws.Range("J31").Value = 3999
$ws.Range("L31").Value = 3999
$ws.Range("N31").Value = -4495
# Row 40
$ws.Range("H40").Value = 3792.75
$ws.Range("I40").Value = 3088.1428
$ws.Range("K40").Value = 3088.1428
$ws.Range("M40").Value = -2952.1428
# Row 82
$ws.Range("H82").Value = 1288.2
$ws.Range("I82").Value = 1169.1428
$ws.Range("K82").Value = 1169.1428
$ws.Range("M82").Value = -808.1428000000001
# Row 85
$ws.Range("H85").Value = 1288.2
$ws.Range("I85").Value = 1169.1428
$ws.Range("K85").Value = 1169.1428
$ws.Range("M85").Value = 78.85719999999992
# Row 93
$ws.Range("H93").Value = 906.5
$ws.Range("J93").Value = 490
$ws.Range("L93").Value = 490
$ws.Range("N93").Value = -2986
# Row 122
$ws.Range("H122").Value = 7932.5483
$ws.Range("I122").Value = 10493.625
$ws.Range("J122").Value = 7041.7393
$ws.Range("K122").Value = 31480.875
$ws.Range("L122").Value = 21125.2179
$ws.Range("M122").Value = -29030.875
$ws.Range("N122").Value = -26025.2179

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 1200
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1200
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1200
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1536
# Row 81
$ws.Range("H81").Value = 2000840.2
$ws.Range("I81").Value = 1049.75
$ws.Range("K81").Value = 2099.5
$ws.Range("M81").Value = -1038.5
# Row 84
$ws.Range("H84").Value = 2000840.2
$ws.Range("I84").Value = 1049.75
$ws.Range("K84").Value = 10497.5
$ws.Range("M84").Value = -5193.5
# Row 96
$ws.Range("H96").Value = 1448.25
$ws.Range("I96").Value = 1200
$ws.Range("J96").Value = 1531
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 1531
$ws.Range("M96").Value = 173
$ws.Range("N96").Value = -4277
# Row 107
$ws.Range("H107").Value = 1744.7333
$ws.Range("I107").Value = 1642.1111
$ws.Range("K107").Value = 4926.3333
$ws.Range("M107").Value = -3006.3333
# Row 122
$ws.Range("H122").Value = 5421.5557
$ws.Range("I122").Value = 3693.2
$ws.Range("K122").Value = 11079.6
$ws.Range("M122").Value = -8629.599999999999
